$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("glory")

# New game (row 14) raw scores
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 23
$ws.Range("C14").Value = 23
$ws.Range("D14").Value = 26
$ws.Range("E14").Value = 11
$ws.Range("F14").Value = 14
$ws.Range("G14").Value = 13

# Running totals
$ws.Range("H14").Formula = "=SUM(C$2:C14)"
$ws.Range("I14").Formula = "=SUM(D$2:D14)"
$ws.Range("J14").Formula = "=SUM(E$2:E14)"
$ws.Range("K14").Formula = "=SUM(F$2:F14)"
$ws.Range("L14").Formula = "=SUM(G$2:G14)"

# Diffs from leader
$ws.Range("M14").Formula = "=H14-MAX(H14:L14)"
$ws.Range("N14").Formula = "=I14-MAX(H14:L14)"
$ws.Range("O14").Formula = "=J14-MAX(H14:L14)"
$ws.Range("P14").Formula = "=K14-MAX(H14:L14)"
$ws.Range("Q14").Formula = "=L14-MAX(H14:L14)"

# Per-game ranks
$ws.Range("R14").Formula = "=RANK(C14,`$C14:`$G14)"
$ws.Range("S14").Formula = "=RANK(D14,`$C14:`$G14)"
$ws.Range("T14").Formula = "=RANK(E14,`$C14:`$G14)"
$ws.Range("U14").Formula = "=RANK(F14,`$C14:`$G14)"
$ws.Range("V14").Formula = "=RANK(G14,`$C14:`$G14)"

# Normalized scores vs glory target
$ws.Range("W14").Formula = "=C14/`$B14"
$ws.Range("X14").Formula = "=D14/`$B14"
$ws.Range("Y14").Formula = "=E14/`$B14"
$ws.Range("Z14").Formula = "=F14/`$B14"
$ws.Range("AA14").Formula = "=G14/`$B14"

$ws.Range("AB14").Formula = "=SUM(C14:G14)"

$ws.Range("X14").Select() | Out-Null
